# "updated to run in local" — refresh the Chrome version under test and
# flip which configured run uses it, then re-point the DATA sheet's
# print orientation and restore the on-screen selections.

$wb = $excel.ActiveWorkbook

$runmanager = $wb.Worksheets.Item("RUNMANAGER")
$data       = $wb.Worksheets.Item("DATA")

# --- DATA sheet: local Chrome version bump + execute-flag swap ---------
# Row 2 (loginLogoutTest / chrome) now runs locally against the newly
# installed Chrome build, so flip it to "yes" and record its version.
$data.Range("B2").Value = "yes"
$data.Range("D2").Value = "'90.0.4430.93"

# The runs that previously rode along on the old chrome 88 build are
# disabled now that row 2 covers it.
$data.Range("B4").Value = "no"
$data.Range("B6").Value = "no"

# DATA sheet should print in portrait orientation.
$data.PageSetup.Orientation = 1

# --- restore on-screen selections left over from the last save ---------
$runmanager.Range("B9").Select() | Out-Null
$data.Range("B7").Select() | Out-Null
